$d = $word.ActiveDocument

# --- First endpoint description (generic path-parameter form) ---
$d.Content.Find.Execute(
    "GET ~/api/v1/Earth/surfaceDistanse/<latA>/<longA>/<latB>/<longB>",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "POST ~/api/v1/Earth/surfaceDistanse^p{^p    “LatitiudeA”: “<latA>”,^p    “longitudeA”: “<longA>”,^p    “LatitiudeB”: “<latB>”,^p    “longitudeB”: “<longB>”^p}^p",
    2) | Out-Null

# --- Second endpoint description (concrete example values) ---
$d.Content.Find.Execute(
    "GET ~/api/v1/ Earth/surfaceDistanse/53.297975/-6.372663/41.385101/-81.440440",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "POST ~/api/v1/ Earth/surfaceDistanse^p{^p    “LatitiudeA”: “53.297975”,^p    “longitudeA”: “-6.372663”,^p    “LatitiudeB”: “41.385101”,^p    “longitudeB”: “-81.440440”^p}",
    2) | Out-Null

Write-Output "done"
